{"js": "// New arithmetic expressions, row-by-row (20 rows x 5 columns), in the\n// same left-to-right / top-to-bottom order as the table cells.\nconst newValues = [\n  [\"19-16=\", \"68-58=\", \"41+33=\", \"86-85=\", \"15+15=\"],\n  [\"23+73=\", \"71-57=\", \"93-75=\", \"38-17=\", \"1+35=\"],\n  [\"75-62=\", \"43+28=\", \"22+14=\", \"14+84=\", \"31+48=\"],\n  [\"9+20=\", \"73-46=\", \"86+4=\", \"69-44=\", \"56-40=\"],\n  [\"61+0=\", \"43-22=\", \"28+15=\", \"93-66=\", \"7+59=\"],\n  [\"15+24=\", \"72-28=\", \"56-23=\", \"16+66=\", \"49+12=\"],\n  [\"15-6=\", \"21+48=\", \"87-67=\", \"74-22=\", \"18+62=\"],\n  [\"59+4=\", \"22+62=\", \"46+0=\", \"32+57=\", \"61+19=\"],\n  [\"63-34=\", \"89-12=\", \"99-65=\", \"69+12=\", \"8+43=\"],\n  [\"45+22=\", \"67-44=\", \"23-6=\", \"44+9=\", \"6+56=\"],\n  [\"53+15=\", \"20+11=\", \"68-50=\", \"85+2=\", \"45+17=\"],\n  [\"79+16=\", \"56-26=\", \"82-14=\", \"36+18=\", \"69-34=\"],\n  [\"9+24=\", \"99-31=\", \"76-27=\", \"83+7=\", \"67-51=\"],\n  [\"67+12=\", \"55+18=\", \"56-56=\", \"19-4=\", \"52+43=\"],\n  [\"57+29=\", \"19-9=\", \"71-25=\", \"78-36=\", \"24+72=\"],\n  [\"28+29=\", \"98-37=\", \"90+2=\", \"8+24=\", \"54-37=\"],\n  [\"49+24=\", \"21+19=\", \"22-2=\", \"24+8=\", \"65+24=\"],\n  [\"13+63=\", \"0+36=\", \"22+50=\", \"44-5=\", \"96-87=\"],\n  [\"5+36=\", \"71-26=\", \"33+45=\", \"25+68=\", \"48+14=\"],\n  [\"92-72=\", \"78-65=\", \"8+10=\", \"6-4=\", \"92-70=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n// Replacing the whole grid in one shot keeps each cell's existing run\n// formatting (font/size) intact while swapping the text content.\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# New arithmetic expressions, row-by-row (20 rows x 5 columns), in the\n# same left-to-right / top-to-bottom order as the table cells.\n$newValues = @(\n    @(\"19-16=\",\"68-58=\",\"41+33=\",\"86-85=\",\"15+15=\"),\n    @(\"23+73=\",\"71-57=\",\"93-75=\",\"38-17=\",\"1+35=\"),\n    @(\"75-62=\",\"43+28=\",\"22+14=\",\"14+84=\",\"31+48=\"),\n    @(\"9+20=\",\"73-46=\",\"86+4=\",\"69-44=\",\"56-40=\"),\n    @(\"61+0=\",\"43-22=\",\"28+15=\",\"93-66=\",\"7+59=\"),\n    @(\"15+24=\",\"72-28=\",\"56-23=\",\"16+66=\",\"49+12=\"),\n    @(\"15-6=\",\"21+48=\",\"87-67=\",\"74-22=\",\"18+62=\"),\n    @(\"59+4=\",\"22+62=\",\"46+0=\",\"32+57=\",\"61+19=\"),\n    @(\"63-34=\",\"89-12=\",\"99-65=\",\"69+12=\",\"8+43=\"),\n    @(\"45+22=\",\"67-44=\",\"23-6=\",\"44+9=\",\"6+56=\"),\n    @(\"53+15=\",\"20+11=\",\"68-50=\",\"85+2=\",\"45+17=\"),\n    @(\"79+16=\",\"56-26=\",\"82-14=\",\"36+18=\",\"69-34=\"),\n    @(\"9+24=\",\"99-31=\",\"76-27=\",\"83+7=\",\"67-51=\"),\n    @(\"67+12=\",\"55+18=\",\"56-56=\",\"19-4=\",\"52+43=\"),\n    @(\"57+29=\",\"19-9=\",\"71-25=\",\"78-36=\",\"24+72=\"),\n    @(\"28+29=\",\"98-37=\",\"90+2=\",\"8+24=\",\"54-37=\"),\n    @(\"49+24=\",\"21+19=\",\"22-2=\",\"24+8=\",\"65+24=\"),\n    @(\"13+63=\",\"0+36=\",\"22+50=\",\"44-5=\",\"96-87=\"),\n    @(\"5+36=\",\"71-26=\",\"33+45=\",\"25+68=\",\"48+14=\"),\n    @(\"92-72=\",\"78-65=\",\"8+10=\",\"6-4=\",\"92-70=\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    $rowValues = $newValues[$r - 1]\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
